$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 343 (C:F) with new value
$ws.Range("C343:F343").Value = 37573500000

# Prep rows 353:355 with the same formatting as row 352 (copies the existing
# cell style, e.g. the date number format on column A) before filling values.
$ws.Range("A352:G352").Copy()
$ws.Range("A353:G355").PasteSpecial(-4122)  # xlPasteFormats

# Append new rows 353-355 with M2 data
$ws.Range("A353").Value = 44986.45833333334
$ws.Range("B353").Value = "ECONOMICS:KWM2"
$ws.Range("C353:F353").Value = 39530200000
$ws.Range("G353").Value = 0

$ws.Range("A354").Value = 45017.45833333334
$ws.Range("B354").Value = "ECONOMICS:KWM2"
$ws.Range("C354:F354").Value = 39461400000
$ws.Range("G354").Value = 0

$ws.Range("A355").Value = 45047.41666666666
$ws.Range("B355").Value = "ECONOMICS:KWM2"
$ws.Range("C355:F355").Value = 39655200000
$ws.Range("G355").Value = 0
